# Update TPM-derived values in the LR-pairs sheet (Bmp8a-Tgfbr1).
# Only the "raw" ligand/receptor expression values actually change;
# all derived specificity / weight columns are recomputed from them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New raw inputs (from updated TPM data)
$newLigandAvg   = 0.061724     # Ligand average expression value (sending cluster = ECs)
$newLigandTotal = 0.185172     # Ligand total expression value   (sending cluster = ECs)
$newReceptorAvg   = 7.106976666666665  # Receptor average expression value (target cluster = ECs)
$newReceptorTotal = 21.32093            # Receptor total expression value   (target cluster = ECs)

# Rows where the sending cluster is "ECs" -> columns G (avg) / H (total) change
$ecsSendingRows = @(2, 3, 4)
foreach ($r in $ecsSendingRows) {
    $ws.Cells.Item($r, 7).Value = $newLigandAvg    # column G
    $ws.Cells.Item($r, 8).Value = $newLigandTotal  # column H
}

# Rows where the target cluster is "ECs" -> columns M (avg) / N (total) change
$ecsTargetRows = @(2, 5, 8)
foreach ($r in $ecsTargetRows) {
    $ws.Cells.Item($r, 13).Value = $newReceptorAvg    # column M
    $ws.Cells.Item($r, 14).Value = $newReceptorTotal  # column N
}

# Recompute the derived columns for every data row (2..10):
#   I = Ligand specificity (avg)      = G_sending / sum(G over all 3 sending clusters)
#   J = Ligand specificity (total)    = H_sending / sum(H over all 3 sending clusters)
#   O = Receptor specificity (avg)    = M_target  / sum(M over all 3 target clusters)
#   P = Receptor specificity (total)  = N_target  / sum(N over all 3 target clusters)
#   Q = Edge average expression weight              = G * M
#   R = Edge total expression weight                = H * N
#   S = Edge average expression derived specificity  = I * O
#   T = Edge total expression derived specificity    = J * P

# Distinct per-sending-cluster ligand values (one representative row per cluster)
$gECs   = $ws.Cells.Item(2, 7).Value2
$gFAPs  = $ws.Cells.Item(5, 7).Value2
$gMuSCs = $ws.Cells.Item(8, 7).Value2
$sumG = $gECs + $gFAPs + $gMuSCs

$hECs   = $ws.Cells.Item(2, 8).Value2
$hFAPs  = $ws.Cells.Item(5, 8).Value2
$hMuSCs = $ws.Cells.Item(8, 8).Value2
$sumH = $hECs + $hFAPs + $hMuSCs

# Distinct per-target-cluster receptor values (one representative row per cluster)
$mECs   = $ws.Cells.Item(2, 13).Value2
$mFAPs  = $ws.Cells.Item(3, 13).Value2
$mMuSCs = $ws.Cells.Item(4, 13).Value2
$sumM = $mECs + $mFAPs + $mMuSCs

$nECs   = $ws.Cells.Item(2, 14).Value2
$nFAPs  = $ws.Cells.Item(3, 14).Value2
$nMuSCs = $ws.Cells.Item(4, 14).Value2
$sumN = $nECs + $nFAPs + $nMuSCs

$sendingClusterByRow = @{ 2="ECs"; 3="ECs"; 4="ECs"; 5="FAPs"; 6="FAPs"; 7="FAPs"; 8="MuSCs"; 9="MuSCs"; 10="MuSCs" }
$targetClusterByRow  = @{ 2="ECs"; 3="FAPs"; 4="MuSCs"; 5="ECs"; 6="FAPs"; 7="MuSCs"; 8="ECs"; 9="FAPs"; 10="MuSCs" }

for ($r = 2; $r -le 10; $r++) {
    $sc = $sendingClusterByRow[$r]
    $tc = $targetClusterByRow[$r]

    if ($sc -eq "ECs") { $g = $gECs; $h = $hECs }
    elseif ($sc -eq "FAPs") { $g = $gFAPs; $h = $hFAPs }
    else { $g = $gMuSCs; $h = $hMuSCs }

    if ($tc -eq "ECs") { $m = $mECs; $n = $nECs }
    elseif ($tc -eq "FAPs") { $m = $mFAPs; $n = $nFAPs }
    else { $m = $mMuSCs; $n = $nMuSCs }

    $i = $g / $sumG
    $j = $h / $sumH
    $o = $m / $sumM
    $p = $n / $sumN

    $ws.Cells.Item($r, 9).Value  = $i            # I
    $ws.Cells.Item($r, 10).Value = $j            # J
    $ws.Cells.Item($r, 15).Value = $o            # O
    $ws.Cells.Item($r, 16).Value = $p            # P
    $ws.Cells.Item($r, 17).Value = $g * $m       # Q
    $ws.Cells.Item($r, 18).Value = $h * $n       # R
    $ws.Cells.Item($r, 19).Value = $i * $o       # S
    $ws.Cells.Item($r, 20).Value = $j * $p       # T
}
